$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forms")

$ws.Range("A5").Value = "AccountDetails"
$ws.Range("B5").Value = "qatesting.lotuswave@gmail.com"
$ws.Range("C5").Value = "Lotuswave@123"
$ws.Range("D5").Value = "Lotuswave@123"
$ws.Range("E5").Value = "QA"
$ws.Range("F5").Value = "TEST"
$ws.Range("H5").Value = "qatesting.lotuswave@gmail.com"

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null

$ws.Range("H12").Select()
